$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the "Folio No" column (old column G). This shifts
#    Instrument/Currency/Investment Domicile (old H/I/J) left to G/H/I.
# ------------------------------------------------------------------
$ws.Range("G1").EntireColumn.Delete()

# ------------------------------------------------------------------
# 2. Fix up the cell comments. Comments stay anchored to their
#    original cell reference when a column is deleted, so:
#      - the stale "Folio No" comment is still sitting on G1
#      - the stale "Investment Domicile" comment is still sitting on J1
#    Move the SEBI/domicile comment text onto the new I1 (Investment
#    Domicile) header and drop the Folio No comment entirely.
# ------------------------------------------------------------------
$domicileCommentText = $ws.Range("J1").Comment.Text()
$ws.Range("J1").Comment.Delete()
$ws.Range("G1").Comment.Delete()
$ws.Range("I1").AddComment($domicileCommentText) | Out-Null

# ------------------------------------------------------------------
# 3. Drop the old sample rows (2-8) completely and start clean.
# ------------------------------------------------------------------
$ws.Range("A2:K8").Clear()
$ws.Range("7:8").EntireRow.Delete()

# ------------------------------------------------------------------
# 4. Re-key the header row labels (order changed after the column
#    delete: Fund | Portfolio Co | Date | Amount | Qty | Notes |
#    Instrument | Currency | Investment Domicile).
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Fund"
$ws.Range("B1").Value = "Portfolio Company Name *"
$ws.Range("C1").Value = "Investment Date *"
$ws.Range("D1").Value = "Amount (Excluding Expenses)*"
$ws.Range("E1").Value = "Quantity *"
$ws.Range("F1").Value = "Notes"
$ws.Range("G1").Value = "Instrument"
$ws.Range("H1").Value = "Currency"
$ws.Range("I1").Value = "Investment Domicile *"

# ------------------------------------------------------------------
# 5. Populate the new sample rows (dates written as serials so the
#    day/month order is unambiguous regardless of runtime locale).
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Demo Fund 1"
$ws.Range("B2").Value = "TSTF1 Port Co 1"
$ws.Range("C2").Value = 45031
$ws.Range("D2").Value = 200000000
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = "Equity"
$ws.Range("H2").Value = "INR"
$ws.Range("I2").Value = "Domestic"

$ws.Range("A3").Value = "Demo Fund 1"
$ws.Range("B3").Value = "TSTF1 Port Co 1"
$ws.Range("C3").Value = 45033
$ws.Range("D3").Value = 150000000
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = "Equity"
$ws.Range("H3").Value = "INR"
$ws.Range("I3").Value = "Domestic"

$ws.Range("A4").Value = "Demo Fund 1"
$ws.Range("B4").Value = "TSTF1 Port Co 1"
$ws.Range("C4").Value = 45031
$ws.Range("D4").Value = 350000000
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = "CCPS"
$ws.Range("H4").Value = "INR"
$ws.Range("I4").Value = "Domestic"

$ws.Range("A5").Value = "Demo Fund 1"
$ws.Range("B5").Value = "TSTF1 Port Co 2"
$ws.Range("C5").Value = 45214
$ws.Range("D5").Value = 1400000000
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = "CCPS"
$ws.Range("H5").Value = "INR"
$ws.Range("I5").Value = "Domestic"

$ws.Range("C2:C5").NumberFormat = "dd/mm/yyyy"

# Quantity = Amount / Price, filled down as a single formula block.
$ws.Range("E2:E5").Formula = "=D2/F2"

# ------------------------------------------------------------------
# 6. Number formatting: Amount (and the helper price cols) use an
#    accounting "Comma" style (thousands separator, 2 decimals).
# ------------------------------------------------------------------
$commaFmt = "_ * #,##0.00_ ;_ * \-#,##0.00_ ;_ * ""-""??_ ;_ @_ "
$ws.Range("D2:D5").NumberFormat = $commaFmt
$ws.Range("F2:F5").NumberFormat = $commaFmt
$ws.Range("J2:K5").NumberFormat = $commaFmt

# ------------------------------------------------------------------
# 7. Update the Investment Domicile data-validation list to track the
#    new column (I) and the row count Excel recalculated it to.
# ------------------------------------------------------------------
$ws.Range("I2:I1008").Validation.Delete()
$ws.Range("I2:I1005").Validation.Add(3, 1, 1, '"Domestic,Overseas"')

# ------------------------------------------------------------------
# 8. Sheet now ends at row 6 (one blank, formatted row under the data).
# ------------------------------------------------------------------
$ws.Range("A6").Value = "Demo Fund 1"
$ws.Range("A6").ClearContents()

$wb.Save()
